$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the last header and add two new trailing headers (PTRES / N. NATUREZA DESPESA / PLANO INTERNO)
$ws.Range("H1").Value = "PTRES"
$ws.Range("I1").Value = "Nº NATUREZA DESPESA"
$ws.Range("J1").Value = "PLANO INTERNO"

# Give the two new header cells the same look (font/fill/border/alignment) as the rest of the header row
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Materialize two blank data rows below the header (rows 2 and 3), mirroring a range that
# already holds genuinely empty/unstyled cells so A2:F3 stay blank with no inherited formatting
$ws.Range("Z1").Copy()
$ws.Range("A2:F3").PasteSpecial(-4122)
